$d = $word.ActiveDocument

# 1. Merge the split "Strengthen customer relationships..." runs back into a
#    single run and drop the mid-word "_GoBack" bookmark that split them.
$d.Content.Find.Execute("reliable delivery of customer", $true, $false, $false, $false, $false, $true, 1, $false, "reliable delivery of customer", 2) | Out-Null

# 2. Append three new paragraphs after the "Transform everyday communications..."
#    paragraph (just before the trailing empty paragraph), all using the
#    Heading4 style inherited from that paragraph.
$transformPara = $d.Paragraphs(5)
$insertionPoint = $transformPara.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()
$insertionPoint.InsertParagraphAfter()
$insertionPoint.InsertParagraphAfter()

# 2a. First new paragraph: empty, but carries the relocated "_GoBack" bookmark.
$bookmarkPara = $d.Paragraphs(6)
$bookmarkRange = $bookmarkPara.Range
$bookmarkRange.Collapse(1)
$bookmarkRange.InsertAfter("X")
$tempRange = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $tempRange)
$delRange = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start + 1)
$delRange.Delete()

# 2b. Second new paragraph: "Additional string" highlighted yellow.
$yellowPara = $d.Paragraphs(7)
$yellowRange = $yellowPara.Range
$yellowRange.Collapse(1)
$yellowRange.InsertAfter("Additional string")
$yellowRange.Font.HighlightColorIndex = 7

# 2c. Third new paragraph: "Plus one additional string" highlighted green.
$greenPara = $d.Paragraphs(8)
$greenRange = $greenPara.Range
$greenRange.Collapse(1)
$greenRange.InsertAfter("Plus one additional string")
$greenRange.Font.HighlightColorIndex = 4
